$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns A, Q, R
$a2 = $ws.Range("A2").Value2
$q2 = $ws.Range("Q2").Value2
$r2 = $ws.Range("R2").Value2

$a3 = $ws.Range("A3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2

$ws.Range("A2").Value2 = $a3
$ws.Range("Q2").Value2 = $q3
$ws.Range("R2").Value2 = $r3

$ws.Range("A3").Value2 = $a2
$ws.Range("Q3").Value2 = $q2
$ws.Range("R3").Value2 = $r2

# Swap values between row 4 and row 5 for columns A, Q, R
$a4 = $ws.Range("A4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

$a5 = $ws.Range("A5").Value2
$q5 = $ws.Range("Q5").Value2
$r5 = $ws.Range("R5").Value2

$ws.Range("A4").Value2 = $a5
$ws.Range("Q4").Value2 = $q5
$ws.Range("R4").Value2 = $r5

$ws.Range("A5").Value2 = $a4
$ws.Range("Q5").Value2 = $q4
$ws.Range("R5").Value2 = $r4
